$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.751.85'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '2.226.64'
$ws.Range("E3").Value = '  -5.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '84.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.513'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.97%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.466'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0787'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -8.49%  '
$ws.Range("E13").Value = '  -2.36%  '
$ws.Range("D14").Value = '2.570.25'
$ws.Range("E14").Value = '  -5.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.77%  '
$ws.Range("D17").Value = '2.222.09'
$ws.Range("E17").Value = '  -6.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.719'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.31%  '
$ws.Range("D19").Value = '39.714.79'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").Value = '0.0₃0884'
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("E21").Value = '  -5.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.13%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.00%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  -5.71%  '
$ws.Range("E35").Value = '  -4.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0701'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.111'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0976'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("E40").Value = '  -5.41%  '
$ws.Range("E41").Value = '  -4.66%  '
$ws.Range("E42").Value = '  -3.69%  '
$ws.Range("D43").Value = '1.951.66'
$ws.Range("E44").Value = '  -3.62%  '
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.00%  '
$ws.Range("E48").Value = '  -4.42%  '
$ws.Range("D49").Value = '2.441.74'
$ws.Range("E49").Value = '  -4.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.82'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.59%  '
